# Generate Report for Handoff
#
# The localization-status report is being regenerated: a new handoff id
# (8638bfc9-2d3f-45f0-a1d6-9fc4dbbfc501) and a new localization-resource
# hash (67025b7aae18e21d6304cb1ae344cc2930ac32a7) replace the previous
# run's values (cade4486-70a4-4ace-b563-62f6f92fdfde /
# aee4ababd66e47d954b7d8d4ff7ea824ff1f91d6), and the handoff timestamps
# move forward slightly. Only the cell text / hyperlink *display* text
# changes - the underlying hyperlink targets (already recorded against
# the previous run's file names) are left exactly as they were.

$oldId   = "cade4486-70a4-4ace-b563-62f6f92fdfde"
$newId   = "8638bfc9-2d3f-45f0-a1d6-9fc4dbbfc501"
$oldHash = "aee4ababd66e47d954b7d8d4ff7ea824ff1f91d6"
$newHash = "67025b7aae18e21d6304cb1ae344cc2930ac32a7"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview" (sheet1): columns File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$overviewLinkAddress = "https://github.com/OpenLocalizationTest/oltest/blob/4d2de3005e7f9bcc52bae4a91a618d7877c06207/e2e/$oldId.md"

$wsOverview.Range("A2").Value = "$newId.md"
$wsOverview.Range("D2").Value = "2016-42-12 08:42:40"

# Rebuild the hyperlink on A2 with the refreshed display text (the
# underlying link target is unchanged). Delete() clears every hyperlink
# on the sheet, so re-add it right away.
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $overviewLinkAddress, "", "", "$newId.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn" (sheet2)
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$zhMdAddress  = "https://github.com/OpenLocalizationTest/oltest/blob/4d2de3005e7f9bcc52bae4a91a618d7877c06207/e2e/$oldId.md"
$zhExtAddress = "https://github.com/OpenLocalizationTest/oltest/blob/4d2de3005e7f9bcc52bae4a91a618d7877c06207/e2e/$oldId.md"
$zhXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dd8d6a0791373ba32c56c27e21a67c0843890c5b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldId.$oldHash.zh-cn.xlf"

$wsZhCn.Range("A2").Value = "$newId.md"
$wsZhCn.Range("D2").Value = "$newId.$newHash.zh-cn.xlf"
$wsZhCn.Range("E2").Value = "2016-03-12 08:42:37"

# Rebuild all three hyperlinks on this sheet (A2, B2, D2); Delete() clears
# the whole collection, so re-add every one in order right after.
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $zhMdAddress, "", "", "$newId.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B2"), $zhExtAddress, "", "", ".md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D2"), $zhXlfAddress, "", "", "$newId.$newHash.zh-cn.xlf")

# ---------------------------------------------------------------------
# Sheet "de-de" (sheet3)
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$deMdAddress  = "https://github.com/OpenLocalizationTest/oltest/blob/4d2de3005e7f9bcc52bae4a91a618d7877c06207/e2e/$oldId.md"
$deExtAddress = "https://github.com/OpenLocalizationTest/oltest/blob/4d2de3005e7f9bcc52bae4a91a618d7877c06207/e2e/$oldId.md"
$deXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/eeb9ce99c8090b1bdc6ec0faab3aa83e566a060e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldId.$oldHash.de-de.xlf"

$wsDeDe.Range("A2").Value = "$newId.md"
$wsDeDe.Range("D2").Value = "$newId.$newHash.de-de.xlf"
$wsDeDe.Range("E2").Value = "2016-03-12 08:42:40"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $deMdAddress, "", "", "$newId.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B2"), $deExtAddress, "", "", ".md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D2"), $deXlfAddress, "", "", "$newId.$newHash.de-de.xlf")

Write-Host "Updated handoff report from $oldId to $newId"
